# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Update the DAMSLTag (column I) and DialogAct (column J) values for the rows
# whose dialog-act annotations changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 19;  I = 'sd';  J = 'Statement-non-opinion' },
    @{ Row = 23;  I = '%';   J = 'Uninterpretable' },
    @{ Row = 28;  I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 49;  I = 'aa';  J = 'Agree/Accept' },
    @{ Row = 62;  I = 'sd';  J = 'Statement-non-opinion' },
    @{ Row = 67;  I = 'sv';  J = 'Statement-opinion' },
    @{ Row = 72;  I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 76;  I = 'ba';  J = 'Appreciation' },
    @{ Row = 83;  I = '%';   J = 'Uninterpretable' },
    @{ Row = 94;  I = 'sd';  J = 'Statement-non-opinion' },
    @{ Row = 106; I = 'sd';  J = 'Statement-non-opinion' },
    @{ Row = 125; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 126; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 127; I = 'ba';  J = 'Appreciation' },
    @{ Row = 132; I = 'aa';  J = 'Agree/Accept' },
    @{ Row = 141; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 146; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 148; I = 'sv';  J = 'Statement-opinion' },
    @{ Row = 182; I = 'sd';  J = 'Statement-non-opinion' },
    @{ Row = 185; I = '%';   J = 'Uninterpretable' },
    @{ Row = 189; I = 'aa';  J = 'Agree/Accept' },
    @{ Row = 195; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 200; I = 'ba';  J = 'Appreciation' },
    @{ Row = 204; I = 'sv';  J = 'Statement-opinion' },
    @{ Row = 206; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 209; I = 'aa';  J = 'Agree/Accept' },
    @{ Row = 215; I = '%';   J = 'Uninterpretable' },
    @{ Row = 220; I = 'sv';  J = 'Statement-opinion' },
    @{ Row = 232; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 236; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 252; I = 'sv';  J = 'Statement-opinion' },
    @{ Row = 259; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 261; I = '%';   J = 'Uninterpretable' },
    @{ Row = 266; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 270; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 285; I = 'b';   J = 'Acknowledge (Backchannel)' },
    @{ Row = 286; I = 'sv';  J = 'Statement-opinion' },
    @{ Row = 288; I = 'sd';  J = 'Statement-non-opinion' },
    @{ Row = 302; I = 'sd';  J = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
